# Auto-committed on 2022/02/24 週四
# Adds a new GenTable entry "SlipEbsRecord" (傳票上傳EBS紀錄檔) to the
# L9-報表作業 section of the table, and refreshes the timestamp on the
# existing TxInquiry (查詢紀錄檔) entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at row 315 (just after the RptSubCom row, right
# before the existing SlipMedia row); this pushes every row from the old
# row 315 onward down by one, preserving their content/styles.
$ws.Rows.Item(315).Insert()

# Fill in the new row with the SlipEbsRecord entry.
$ws.Cells.Item(315, 1).Value = "L9-報表作業"
$ws.Cells.Item(315, 2).Value = "SlipEbsRecord"
$ws.Cells.Item(315, 3).Value = "傳票上傳EBS紀錄檔"
$ws.Cells.Item(315, 4).Formula = '=HYPERLINK("[\\192.168.10.16\St1Share(NAS)\SKL\DB\GenTables\L9-報表作業\SlipEbsRecord.xlsx]DBD!A1", "連結")'
$ws.Cells.Item(315, 5).Value = "2022年02月24日 11:10:20"

# The TxInquiry (查詢紀錄檔) row, originally row 340, is now row 341 after
# the insertion above; update its "last generated" timestamp.
$ws.Cells.Item(341, 5).Value = "2022年02月24日 10:25:59"
